$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "From" bound of rule R30 (cell C10) from 18 to 1
$ws.Range("C10").Value = 1.0
